$wb = $excel.ActiveWorkbook

# The "想去人数" (interest count) values in column F changed for both the
# "展览" and "全部类型" sheets (they mirror the same data).
$updates = @{
    2  = 369
    3  = 365
    4  = 1878
    5  = 80
    7  = 193
    8  = 747
    11 = 4470
    14 = 1249
    17 = 835
    19 = 455
    21 = 223
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
